$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-PlainCell {
    param($addr, $val)
    $ws.Range($addr).Value = $val
}

# Row 2
Set-TextCell "D2" "60.866.69"

# Row 3
Set-TextCell "D3" "2.401.18"
Set-PlainCell "E3" "  -0.40%  "

# Row 4
Set-PlainCell "E4" "  +0.49%  "

# Row 5
Set-TextCell "D5" "566.20"

# Row 6
Set-TextCell "D6" "141.97"
Set-PlainCell "E6" "  +2.74%  "

# Row 7
Set-PlainCell "E7" "  -0.37%  "

# Row 8
Set-PlainCell "E8" "  +2.70%  "

# Row 9
Set-TextCell "D9" "2.408.45"
Set-PlainCell "E9" "  +0.59%  "

# Row 10
Set-PlainCell "E10" "  +2.18%  "

# Row 11
Set-PlainCell "E11" "  -0.05%  "

# Row 12
Set-PlainCell "E12" "  +2.96%  "

# Row 13
Set-PlainCell "E13" "  +2.81%  "

# Row 14
Set-TextCell "D14" "26.46"
Set-PlainCell "E14" "  +2.35%  "

# Row 15
Set-TextCell "D15" "0.0000170"
Set-PlainCell "E15" "  +0.24%  "

# Row 16
Set-TextCell "D16" "2.835.80"
Set-PlainCell "E16" "  -0.73%  "

# Row 17
Set-TextCell "D17" "60.781.08"
Set-PlainCell "E17" "  +0.01%  "

# Row 18
Set-TextCell "D18" "2.408.17"
Set-PlainCell "E18" "  +0.50%  "

# Row 19
Set-TextCell "D19" "8.06"
Set-PlainCell "E19" "  +3.52%  "

# Row 20
Set-TextCell "D20" "10.70"

# Row 21
Set-TextCell "D21" "324.26"
Set-PlainCell "E21" "  +0.78%  "

# Row 22
Set-PlainCell "E22" "  +1.79%  "

# Row 23
Set-TextCell "D23" "6.03"
Set-PlainCell "E23" "  -0.98%  "

# Row 24
Set-PlainCell "E24" "  -0.18%  "

# Row 25
Set-PlainCell "E25" "  +5.04%  "

# Row 26
Set-TextCell "D26" "65.10"
Set-PlainCell "E26" "  +0.30%  "

# Row 27
Set-TextCell "D27" "584.25"
Set-PlainCell "E27" "  +1.53%  "

# Row 28
Set-TextCell "D28" "8.23"
Set-PlainCell "E28" "  +0.55%  "

# Row 29
Set-PlainCell "B29" "WrappedeETH"
Set-PlainCell "C29" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextCell "D29" "2.518.46"
Set-PlainCell "E29" "  -0.70%  "

# Row 30
Set-PlainCell "B30" "PEPE"
Set-PlainCell "C30" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D30" "0.0₃0943"
Set-PlainCell "E30" "  +3.26%  "

# Row 31
Set-TextCell "D31" "8.03"
Set-PlainCell "E31" "  +3.28%  "

# Row 32
Set-PlainCell "E32" "  +1.69%  "

# Row 33
Set-TextCell "D33" "1.82"
Set-PlainCell "E33" "  +0.33%  "

# Row 34
Set-PlainCell "E34" "  +0.85%  "

# Row 35
Set-PlainCell "E35" "  +6.02%  "

# Row 36
Set-PlainCell "E36" "  -0.72%  "

# Row 37
Set-TextCell "D37" "153.64"
Set-PlainCell "E37" "  +1.04%  "

# Row 38
Set-TextCell "D38" "0.372"
Set-PlainCell "E38" "  +1.97%  "

# Row 39
Set-PlainCell "E39" "  +1.56%  "

# Row 40
Set-TextCell "D40" "18.31"
Set-PlainCell "E40" "  +0.87%  "

# Row 41
Set-TextCell "D41" "5.20"
Set-PlainCell "E41" "  +2.33%  "

# Row 43
Set-PlainCell "B43" "Stacks"
Set-PlainCell "C43" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D43" "1.69"
Set-PlainCell "E43" "  +2.25%  "

# Row 44
Set-PlainCell "B44" "dogwifhat"
Set-PlainCell "C44" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D44" "2.52"
Set-PlainCell "E44" "  +12.03%  "

# Row 45
Set-TextCell "D45" "41.64"
Set-PlainCell "E45" "  +1.35%  "

# Row 46
Set-TextCell "D46" "0.0₆0282"
Set-PlainCell "E46" "  +8.08%  "

# Row 47
Set-TextCell "D47" "141.61"
Set-PlainCell "E47" "  +0.22%  "

# Row 48
Set-PlainCell "E48" "  +1.07%  "

# Row 49
Set-TextCell "D49" "0.592"
Set-PlainCell "E49" "  +1.46%  "

# Row 50
Set-TextCell "D50" "0.0510"
Set-PlainCell "E50" "  +2.51%  "

# Row 51
Set-TextCell "D51" "19.55"
Set-PlainCell "E51" "  +1.95%  "
